$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the weekly data block (row 969),
# pushing all existing data (previously rows 969:1046) down to 971:1048.
$ws.Range("A969:A970").EntireRow.Insert()

# New week's data (Primera / Segunda) for the Brocoli price table.
# Row 969 - Primera
$ws.Cells.Item(969, 1).Value = 8
$ws.Cells.Item(969, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(969, 3).Value = "Coquimbo"
$ws.Cells.Item(969, 4).Value = 45013
$ws.Cells.Item(969, 5).Value = 4
$ws.Cells.Item(969, 6).Value = 100112023
$ws.Cells.Item(969, 7).Value = "Brócoli"
$ws.Cells.Item(969, 8).Value = "Sin especificar"
$ws.Cells.Item(969, 9).Value = "Primera"
$ws.Cells.Item(969, 10).Value = 2100
$ws.Cells.Item(969, 11).Value = 900
$ws.Cells.Item(969, 12).Value = 1000
$ws.Cells.Item(969, 13).Value = 950
$ws.Cells.Item(969, 14).Value = "$/unidad"
$ws.Cells.Item(969, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(969, 16).Value = 950
$ws.Cells.Item(969, 17).Value = 1
$ws.Cells.Item(969, 18).Value = "Hortaliza"

# Row 970 - Segunda
$ws.Cells.Item(970, 1).Value = 8
$ws.Cells.Item(970, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(970, 3).Value = "Coquimbo"
$ws.Cells.Item(970, 4).Value = 45013
$ws.Cells.Item(970, 5).Value = 4
$ws.Cells.Item(970, 6).Value = 100112023
$ws.Cells.Item(970, 7).Value = "Brócoli"
$ws.Cells.Item(970, 8).Value = "Sin especificar"
$ws.Cells.Item(970, 9).Value = "Segunda"
$ws.Cells.Item(970, 10).Value = 1480
$ws.Cells.Item(970, 11).Value = 700
$ws.Cells.Item(970, 12).Value = 800
$ws.Cells.Item(970, 13).Value = 750
$ws.Cells.Item(970, 14).Value = "$/unidad"
$ws.Cells.Item(970, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(970, 16).Value = 750
$ws.Cells.Item(970, 17).Value = 1
$ws.Cells.Item(970, 18).Value = "Hortaliza"

# Make sure the date cells keep the date number format used elsewhere in column D.
$ws.Range("D969:D970").NumberFormat = $ws.Range("D971").NumberFormat
